$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114
$ws.Cells.Item(114, 1).Value = 112
$ws.Cells.Item(114, 2).Value = 'N/A'
$ws.Cells.Item(114, 3).Value = 'N/A'
$ws.Cells.Item(114, 4).Value = 'No'
$ws.Cells.Item(114, 5).Value = 'No'
$ws.Cells.Item(114, 6).Value = 'No'
$ws.Cells.Item(114, 7).Value = 'No'
$ws.Cells.Item(114, 8).Value = 'N/A'
$ws.Cells.Item(114, 9).Value = 'No'
$ws.Cells.Item(114, 10).Value = 'Day'
$ws.Cells.Item(114, 11).Value = 'No'
$ws.Cells.Item(114, 12).Value = 'No accident occurred'

# Row 115
$ws.Cells.Item(115, 1).Value = 113
$ws.Cells.Item(115, 2).Value = 2
$ws.Cells.Item(115, 3).Value = 'rear-end'
$ws.Cells.Item(115, 4).Value = 'No'
$ws.Cells.Item(115, 5).Value = 'No'
$ws.Cells.Item(115, 6).Value = 'No'
$ws.Cells.Item(115, 7).Value = 'No'
$ws.Cells.Item(115, 8).Value = 'sedan, suv'
$ws.Cells.Item(115, 9).Value = 'No'
$ws.Cells.Item(115, 10).Value = 'Night'
$ws.Cells.Item(115, 11).Value = 'No'
$ws.Cells.Item(115, 12).Value = 'Fender Bender'

# Row 116
$ws.Cells.Item(116, 1).Value = 114
$ws.Cells.Item(116, 2).Value = 2
$ws.Cells.Item(116, 3).Value = 'rear-end'
$ws.Cells.Item(116, 4).Value = 'No'
$ws.Cells.Item(116, 5).Value = 'No'
$ws.Cells.Item(116, 6).Value = 'No'
$ws.Cells.Item(116, 7).Value = 'No'
$ws.Cells.Item(116, 8).Value = 'sedan, suv'
$ws.Cells.Item(116, 9).Value = 'No'
$ws.Cells.Item(116, 10).Value = 'Night'
$ws.Cells.Item(116, 11).Value = 'No'
$ws.Cells.Item(116, 12).Value = 'Fender Bender'

# Row 117
$ws.Cells.Item(117, 1).Value = 115
$ws.Cells.Item(117, 2).Value = 2
$ws.Cells.Item(117, 3).Value = 'slight t-bone'
$ws.Cells.Item(117, 4).Value = 'No'
$ws.Cells.Item(117, 5).Value = 'No'
$ws.Cells.Item(117, 6).Value = 'No'
$ws.Cells.Item(117, 7).Value = 'Yes'
$ws.Cells.Item(117, 8).Value = 'sedan, van'
$ws.Cells.Item(117, 9).Value = 'No'
$ws.Cells.Item(117, 10).Value = 'Day'
$ws.Cells.Item(117, 11).Value = 'No'
$ws.Cells.Item(117, 12).Value = 'Roundabout'

# Row 118
$ws.Cells.Item(118, 1).Value = 116
$ws.Cells.Item(118, 2).Value = 2
$ws.Cells.Item(118, 3).Value = 't-bone'
$ws.Cells.Item(118, 4).Value = 'Yes'
$ws.Cells.Item(118, 5).Value = 'Yes'
$ws.Cells.Item(118, 6).Value = 'No'
$ws.Cells.Item(118, 7).Value = 'Yes'
$ws.Cells.Item(118, 8).Value = 'van x2'
$ws.Cells.Item(118, 9).Value = 'No'
$ws.Cells.Item(118, 10).Value = 'Day'
$ws.Cells.Item(118, 11).Value = 'No'
$ws.Cells.Item(118, 12).Value = 'intersection'

# Row 119
$ws.Cells.Item(119, 1).Value = 117
$ws.Cells.Item(119, 2).Value = 2
$ws.Cells.Item(119, 3).Value = 'slight t-bone'
$ws.Cells.Item(119, 4).Value = 'No'
$ws.Cells.Item(119, 5).Value = 'No'
$ws.Cells.Item(119, 6).Value = 'No'
$ws.Cells.Item(119, 7).Value = 'Yes'
$ws.Cells.Item(119, 8).Value = 'sedan x2'
$ws.Cells.Item(119, 9).Value = 'No'
$ws.Cells.Item(119, 10).Value = 'Day'
$ws.Cells.Item(119, 11).Value = 'No'
$ws.Cells.Item(119, 12).Value = 'intersection'

# Row 120
$ws.Cells.Item(120, 1).Value = 118
$ws.Cells.Item(120, 2).Value = 2
$ws.Cells.Item(120, 3).Value = 'rear-end'
$ws.Cells.Item(120, 4).Value = 'Yes'
$ws.Cells.Item(120, 5).Value = 'Yes'
$ws.Cells.Item(120, 6).Value = 'No'
$ws.Cells.Item(120, 7).Value = 'Yes'
$ws.Cells.Item(120, 8).Value = 'motorbike, suv'
$ws.Cells.Item(120, 9).Value = 'No'
$ws.Cells.Item(120, 10).Value = 'Day'
$ws.Cells.Item(120, 11).Value = 'No'
$ws.Cells.Item(120, 12).Value = 'Motorcyclist runs into suv and falls off'

# Row 121
$ws.Cells.Item(121, 1).Value = 119
$ws.Cells.Item(121, 2).Value = 2
$ws.Cells.Item(121, 3).Value = 't-bone'
$ws.Cells.Item(121, 4).Value = 'Yes'
$ws.Cells.Item(121, 5).Value = 'Yes'
$ws.Cells.Item(121, 6).Value = 'No'
$ws.Cells.Item(121, 7).Value = 'Yes'
$ws.Cells.Item(121, 8).Value = 'sedan x2'
$ws.Cells.Item(121, 9).Value = 'No'
$ws.Cells.Item(121, 10).Value = 'Day, rainy'
$ws.Cells.Item(121, 11).Value = 'No'
# NOTE: L122 is written here (ahead of L121) so the new shared strings get
# allocated in the same order as the source workbook (row 122's "Other"
# text was authored before row 121's in the original edit history).
$ws.Cells.Item(122, 12).Value = 'intsection at a light'
$ws.Cells.Item(121, 12).Value = 'intersection/roundabout'

# Row 122
$ws.Cells.Item(122, 1).Value = 120
$ws.Cells.Item(122, 2).Value = 2
$ws.Cells.Item(122, 3).Value = 't-bone'
$ws.Cells.Item(122, 4).Value = 'Yes'
$ws.Cells.Item(122, 5).Value = 'Yes'
$ws.Cells.Item(122, 6).Value = 'No'
$ws.Cells.Item(122, 7).Value = 'Yes'
$ws.Cells.Item(122, 8).Value = 'van x2'
$ws.Cells.Item(122, 9).Value = 'No'
$ws.Cells.Item(122, 10).Value = 'Day, rainy'
$ws.Cells.Item(122, 11).Value = 'No'

# Row 123
$ws.Cells.Item(123, 1).Value = 121
$ws.Cells.Item(123, 2).Value = 2
$ws.Cells.Item(123, 3).Value = 't-bone'
$ws.Cells.Item(123, 4).Value = 'Yes'
$ws.Cells.Item(123, 5).Value = 'Yes'
$ws.Cells.Item(123, 6).Value = 'No'
$ws.Cells.Item(123, 7).Value = 'Yes'
$ws.Cells.Item(123, 8).Value = 'suv, sedan'
$ws.Cells.Item(123, 9).Value = 'No'
$ws.Cells.Item(123, 10).Value = 'Day, clear'
$ws.Cells.Item(123, 11).Value = 'No'
$ws.Cells.Item(123, 12).Value = 'intersection'

# Row 124
$ws.Cells.Item(124, 1).Value = 122
$ws.Cells.Item(124, 2).Value = 2
$ws.Cells.Item(124, 3).Value = 't-bone'
$ws.Cells.Item(124, 4).Value = 'Yes'
$ws.Cells.Item(124, 5).Value = 'Yes'
$ws.Cells.Item(124, 6).Value = 'No'
$ws.Cells.Item(124, 7).Value = 'Yes'
$ws.Cells.Item(124, 8).Value = 'sedan x2'
$ws.Cells.Item(124, 9).Value = 'No'
$ws.Cells.Item(124, 10).Value = 'Night'
$ws.Cells.Item(124, 11).Value = 'Yes'
$ws.Cells.Item(124, 12).Value = 'intersection'

# Row 125
$ws.Cells.Item(125, 1).Value = 123
$ws.Cells.Item(125, 2).Value = 2
$ws.Cells.Item(125, 3).Value = 'light side-on-side collision'
$ws.Cells.Item(125, 4).Value = 'No'
$ws.Cells.Item(125, 5).Value = 'No'
$ws.Cells.Item(125, 6).Value = 'No'
$ws.Cells.Item(125, 7).Value = 'Yes'
$ws.Cells.Item(125, 8).Value = 'sedan, bus'
$ws.Cells.Item(125, 9).Value = 'No'
$ws.Cells.Item(125, 10).Value = 'Day'
$ws.Cells.Item(125, 11).Value = 'No'
$ws.Cells.Item(125, 12).Value = 'light tap'

# Row 126
$ws.Cells.Item(126, 1).Value = 124
$ws.Cells.Item(126, 2).Value = 2
$ws.Cells.Item(126, 3).Value = 't-bone'
$ws.Cells.Item(126, 4).Value = 'Yes'
$ws.Cells.Item(126, 5).Value = 'Yes'
$ws.Cells.Item(126, 6).Value = 'No'
$ws.Cells.Item(126, 7).Value = 'Yes'
$ws.Cells.Item(126, 8).Value = 'suv x2'
$ws.Cells.Item(126, 9).Value = 'No'
$ws.Cells.Item(126, 10).Value = 'Day, cloudy'
$ws.Cells.Item(126, 11).Value = 'No'
$ws.Cells.Item(126, 12).Value = 'hit car spins out'

# Row 127
$ws.Cells.Item(127, 1).Value = 125
$ws.Cells.Item(127, 2).Value = 2
$ws.Cells.Item(127, 3).Value = 'front-end collision'
$ws.Cells.Item(127, 4).Value = 'Yes'
$ws.Cells.Item(127, 5).Value = 'Yes'
$ws.Cells.Item(127, 6).Value = 'No'
$ws.Cells.Item(127, 7).Value = 'Yes'
$ws.Cells.Item(127, 8).Value = 'truck, bike'
$ws.Cells.Item(127, 9).Value = 'No'
$ws.Cells.Item(127, 10).Value = 'Day, cloudy'
$ws.Cells.Item(127, 11).Value = 'No'
$ws.Cells.Item(127, 12).Value = 'tuck hits and topples bike'

# Update selection to match the final workbook state
$null = $ws.Range("A128:A134").Select()

Write-Output "done"